# edit.ps1 - apply the "Updated the schema info and diagram to match what we
# have" commit to the Relational Schema document via Word COM interop.
#
# Summary of semantic edits performed:
#   1. RequestQueue(...) second field: "SongID" -> "VersionID"
#   2. SongContributors(...): remove the (single) underline from "Role"
#   3. Files paragraph: "a Source to locate the file" -> "a FileName to locate the file"
#   4. Queue paragraph: ", or dj next queue." -> ", now playing queue or history queue."
#   5. Two cosmetic run-merges (no text change) that happen to fall out of the
#      same editing session, reproduced so the underlying runs collapse the
#      same way Word's own editor would collapse them:
#        - "...Description)" + " "  -> single run "...Description) "
#        - "†" + " " -> single run "† "

$d = $word.ActiveDocument

function Find-Unique {
    <#
        Returns a Range for the first match of $text in the document,
        scanning forward from the very start of the story. Throws if no
        match is found, so typos fail loudly instead of silently no-op'ing.
    #>
    param(
        [string]$text
    )
    $rng = $d.Content
    $rng.Start = 0
    $rng.End = $d.Content.End
    $found = $rng.Find.Execute($text, $false, $false, $false, $false, $false, $true, 0, $false, "", 0)
    if (-not $found) {
        throw "Find-Unique: could not find '$text'"
    }
    return $rng
}

function Replace-Text {
    <#
        Plain text -> text replacement. Leaves Word free to choose how the
        backing runs are split/merged, exactly like typing a replacement
        into a single contiguous Find match would.
    #>
    param(
        [string]$old,
        [string]$new
    )
    $rng = Find-Unique $old
    $rng.Text = $new
}

function Merge-Runs {
    <#
        Forces Word to collapse the run(s) spanning $text into a single run
        by round-tripping the content through a throwaway placeholder. This
        mirrors what the real editor does when you retype text that used to
        straddle a run boundary (e.g. after toggling formatting back off, or
        after accepting an autocorrect) - the adjacent runs with identical
        formatting re-coalesce into one <w:r>.
    #>
    param(
        [string]$text
    )
    $rng = Find-Unique $text
    $rng.Text = "ZZ__MERGE_PLACEHOLDER__ZZ"
    $rng2 = Find-Unique "ZZ__MERGE_PLACEHOLDER__ZZ"
    $rng2.Text = $text
}

# ---------------------------------------------------------------------------
# 1. RequestQueue(RequestID, SongID..., UserID..., ...) -> VersionID
#    (Only the RequestQueue occurrence - the SongVersions/SongContributors/
#    Songs/"attach the SongID" occurrences of "SongID" must stay untouched.)
# ---------------------------------------------------------------------------
$rng = Find-Unique "RequestID, SongID"
$rng.Start = $rng.End - 6   # narrow down to just the "SongID" tail
$rng.Text = "VersionID"

# ---------------------------------------------------------------------------
# 2. SongContributors(..., Role) - drop the underline direct formatting on
#    "Role" (it now matches the other attribute names, which are plain).
# ---------------------------------------------------------------------------
$rng = Find-Unique ", Role)"
$roleStart = $rng.Start + 2      # skip ", "
$roleEnd = $rng.End - 1          # drop trailing ")"
$roleRng = $d.Range($roleStart, $roleEnd)
$roleRng.Font.Underline = 0      # wdUnderlineNone

# ---------------------------------------------------------------------------
# 3. "...hold the FileID, a Source to locate the file, and a Description..."
#    -> "...a FileName to locate the file..."
# ---------------------------------------------------------------------------
Replace-Text "a Source to locate the file" "a FileName to locate the file"

# ---------------------------------------------------------------------------
# 4. "...priority queue, or dj next queue. The AmountPaid..."
#    -> "...priority queue, now playing queue or history queue. The AmountPaid..."
# ---------------------------------------------------------------------------
Replace-Text ", or dj next queue. The" ", now playing queue or history queue. The"

# ---------------------------------------------------------------------------
# 5. Cosmetic run merges (text unchanged either side).
# ---------------------------------------------------------------------------
Merge-Runs ", Description) "

$dagger = [char]0x2020
Merge-Runs "$dagger "
